$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new profit row (row 82) for the 2025-11-07 run.
# Column A holds the date as literal text (matches the existing rows, which
# are plain strings rather than real date cells), so we momentarily force a
# text format before assigning the value to stop Excel's autocorrect from
# converting "11/07/2025" into a date serial, then clear the formatting back
# off again so the cell ends up with the same default style as its neighbors.
$ws.Range("A82").NumberFormat = "@"
$ws.Range("A82").Value = "11/07/2025"
$ws.Range("A82").ClearFormats()

$ws.Range("B82").Value = 9602.040000000001
